$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the current column G (cl-p12-au etc. shift right to I:M)
$ws.Range("G1:H1").EntireColumn.Insert()

# --- New header cells for cl_x2c_acv3z_1s / cl_x2c_acv3z_2s (columns G, H) ---
$ws.Range("G1").Value = "cl_x2c_acv3z_1s"
$ws.Range("H1").Value = "cl_x2c_acv3z_2s"
$ws.Range("H1").Font.Color = 0

# --- New header cells for the HCl block (columns N..W) ---
$ws.Range("N1").Value = "zero_hcl_zora_1s"
$ws.Range("O1").Value = "zero_hcl_zora_2s"
$ws.Range("P1").Value = "zero_hcl_zora_2p12"
$ws.Range("Q1").Value = "zero_hcl_zora_2p32"
$ws.Range("R1").Value = "zero_hcl_zora_2p"
$ws.Range("S1").Value = "hcl_x2c_acv3z_1s"
$ws.Range("T1").Value = "hcl_x2c_acv3z_2s"
$ws.Range("T1").Font.Color = 0
$ws.Range("U1").Value = "hcl_x2c_acv3z_2p_1half"
$ws.Range("V1").Value = "hcl_x2c_acv3z_2p_3half"
$ws.Range("W1").Value = "hcl_x2c_acv3z_2p"

# --- Data values for new columns G,H (rows 2-7) ---
$ws.Range("G2").Value = 2824.94
$ws.Range("H2").Value = 270.33999999999997
$ws.Range("G3").Value = 2825.73
$ws.Range("H3").Value = 271.13
$ws.Range("G4").Value = 2825.92
$ws.Range("H4").Value = 271.32
$ws.Range("G5").Value = 2826.13
$ws.Range("H5").Value = 271.54000000000002
$ws.Range("G6").Value = 2825.92
$ws.Range("H6").Value = 271.32
$ws.Range("G7").Value = 2825.94
$ws.Range("H7").Value = 271.38

# --- Data values for new HCl columns N..V and formula column W (rows 2-7) ---
$ws.Range("N2").Value = 2764.75
$ws.Range("O2").Value = 254.19
$ws.Range("P2").Value = 194.07
$ws.Range("Q2").Value = 192.38
$ws.Range("R2").Value = 193.23
$ws.Range("S2").Value = 2834.68
$ws.Range("T2").Value = 280.3
$ws.Range("U2").Value = 209.74
$ws.Range("V2").Value = 208.08
$ws.Range("W2").Formula = "=(U2+V2)/2"

$ws.Range("N3").Value = 2764.75
$ws.Range("O3").Value = 254.09
$ws.Range("P3").Value = 194.26
$ws.Range("Q3").Value = 192.56
$ws.Range("R3").Value = 193.41
$ws.Range("S3").Value = 2834.81
$ws.Range("T3").Value = 280.45999999999998
$ws.Range("U3").Value = 209.52
$ws.Range("V3").Value = 207.86
$ws.Range("W3").Formula = "=(U3+V3)/2"

$ws.Range("N4").Value = 2764.75
$ws.Range("O4").Value = 254.12
$ws.Range("P4").Value = 194.29
$ws.Range("Q4").Value = 192.57
$ws.Range("R4").Value = 193.43
$ws.Range("S4").Value = 2834.84
$ws.Range("T4").Value = 280.49
$ws.Range("U4").Value = 209.55
$ws.Range("V4").Value = 207.89
$ws.Range("W4").Formula = "=(U4+V4)/2"

$ws.Range("N5").Value = 2764.75
$ws.Range("O5").Value = 254.16
$ws.Range("P5").Value = 194.32
$ws.Range("Q5").Value = 192.62
$ws.Range("R5").Value = 193.47
$ws.Range("S5").Value = 2834.84
$ws.Range("T5").Value = 280.49
$ws.Range("U5").Value = 209.55
$ws.Range("V5").Value = 207.89
$ws.Range("W5").Formula = "=(U5+V5)/2"

$ws.Range("N6").Value = 2764.45
$ws.Range("O6").Value = 253.93
$ws.Range("P6").Value = 194.29
$ws.Range("Q6").Value = 192.59
$ws.Range("R6").Value = 193.44
$ws.Range("S6").Value = 2834.62
$ws.Range("T6").Value = 280.27
$ws.Range("U6").Value = 209.36
$ws.Range("V6").Value = 207.7
$ws.Range("W6").Formula = "=(U6+V6)/2"

$ws.Range("N7").Value = 2764.45
$ws.Range("O7").Value = 253.98
$ws.Range("P7").Value = 194.31
$ws.Range("Q7").Value = 192.61
$ws.Range("R7").Value = 193.46
$ws.Range("S7").Value = 2834.65
$ws.Range("T7").Value = 280.3
$ws.Range("U7").Value = 209.39
$ws.Range("V7").Value = 207.72
$ws.Range("W7").Formula = "=(U7+V7)/2"

# --- sheet view / selection to match target ---
# (scroll the window so column P is near the left edge, then leave the
#  selection on X9, matching the saved view state of the edited workbook)
$excel.ActiveWindow.ScrollColumn = 16
$ws.Range("X9").Select()
